$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Battle")
$rng = $ws.Range("K2:N5")
$rng.BorderAround(1, 2)  # xlContinuous=1, xlThin=2
Write-Host "done"
